$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new "TuVi" sheet right after Sheet1 and make it the active sheet
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$newSheet.Name = "TuVi"

$newSheet.Range("A1").Value = 'Con_Giap'
$newSheet.Range("B1").Value = 'Tong_Quan'
$newSheet.Range("C1").Value = 'Su_Nghiep'
$newSheet.Range("D1").Value = 'Tai_Loc'
$newSheet.Range("E1").Value = 'Tinh_Cam'

$newSheet.Range("A2").Value = 'Tý'
$newSheet.Range("B2").Value = 'Năm 2026 Bính Ngọ là năm Xung Thái Tuế với tuổi Tý. Cuộc sống có nhiều biến động lớn, "thay da đổi thịt". Cần giữ tâm vững vàng trước sóng gió.'
$newSheet.Range("C2").Value = 'Công việc gặp nhiều trắc trở, dễ bị tiểu nhân quấy phá hoặc thay đổi vị trí công tác bất ngờ. "Án binh bất động" là thượng sách.'
$newSheet.Range("D2").Value = 'Tiền bạc vào cửa trước ra cửa sau. Tránh đầu tư mạo hiểm, đặc biệt là chứng khoán hay bất động sản lướt sóng.'
$newSheet.Range("E2").Value = 'Dễ có mâu thuẫn, cãi vã do áp lực cuộc sống. Cần nhường nhịn để giữ hòa khí gia đình.'

$newSheet.Range("A3").Value = 'Sửu'
$newSheet.Range("B3").Value = 'Năm Tương Hại Thái Tuế. Tuy không quá xấu như Tý nhưng dễ gặp chuyện thị phi, "làm ơn mắc oán".'
$newSheet.Range("C3").Value = 'Công việc trì trệ, nỗ lực nhiều nhưng kết quả chưa xứng đáng. Cẩn trọng giấy tờ, ký tá hợp đồng.'
$newSheet.Range("D3").Value = 'Tài lộc trung bình. Chỉ nên tập trung vào công việc chính, chớ ôm đồm nghề tay trái kẻo "xôi hỏng bỏng không".'
$newSheet.Range("E3").Value = 'Tình cảm có phần lạnh nhạt. Người độc thân kén chọn, khó tìm được ý trung nhân ưng ý.'

$newSheet.Range("A4").Value = 'Dần'
$newSheet.Range("B4").Value = 'Năm Tam Hợp (Dần - Ngọ - Tuất). Đây là năm Đại Cát để bứt phá. Thiên thời - Địa lợi - Nhân hòa đều hội tụ.'
$newSheet.Range("C4").Value = 'Có quý nhân phù trợ, thăng quan tiến chức hoặc mở rộng quy mô kinh doanh cực tốt. Làm đâu thắng đó.'
$newSheet.Range("D4").Value = 'Tiền bạc dồi dào, có lộc từ phương xa đưa tới hoặc trúng thưởng, thừa kế.'
$newSheet.Range("E4").Value = 'Đỏ tình đỏ cả bạc. Gia đạo êm ấm, có hỷ tín cưới xin hoặc đón thêm thành viên mới.'

$newSheet.Range("A5").Value = 'Mão'
$newSheet.Range("B5").Value = 'Năm Phá Thái Tuế. Cẩn thận các mối quan hệ xã giao, dễ bị bạn bè lừa gạt hoặc đâm sau lưng.'
$newSheet.Range("C5").Value = 'Công việc có sự thay đổi (chuyển việc, đổi chỗ ngồi). Cần khiêm tốn, tránh phô trương kẻo bị ghen ghét.'
$newSheet.Range("D5").Value = 'Hao tài tốn của vì những chuyện không đâu (sửa xe, thuốc men, tiệc tùng). Cần lập kế hoạch chi tiêu chặt chẽ.'
$newSheet.Range("E5").Value = 'Dễ nảy sinh nghi ngờ, ghen tuông vô cớ. Cần tin tưởng đối phương hơn.'

$newSheet.Range("A6").Value = 'Thìn'
$newSheet.Range("B6").Value = 'Một năm Bình Hòa, mọi thứ ở mức ổn định. Không quá đột phá nhưng cũng không có biến cố lớn. Là khoảng lặng để tích lũy.'
$newSheet.Range("C6").Value = 'Công việc diễn ra đều đều. Nếu muốn khởi nghiệp thì cần chuẩn bị kỹ, chưa phải thời điểm vàng để bung lụa.'
$newSheet.Range("D6").Value = 'Thu nhập ổn định từ lương cứng. Không có nhiều khoản lộc bất ngờ.'
$newSheet.Range("E6").Value = 'Tình cảm êm đềm. Là năm tốt để hâm nóng tình cảm vợ chồng bằng những chuyến du lịch ngắn.'

$newSheet.Range("A7").Value = 'Tỵ'
$newSheet.Range("B7").Value = 'Năm bản lề trước khi bước vào năm tuổi. Có nhiều cơ hội mới mở ra nhưng áp lực công việc cực lớn.'
$newSheet.Range("C7").Value = 'Được cấp trên tin tưởng giao trọng trách. Tuy vất vả nhưng học hỏi được nhiều kinh nghiệm quý báu.'
$newSheet.Range("D7").Value = 'Tiền kiếm được nhiều nhưng chi tiêu cũng lắm (mua sắm tài sản lớn, sửa nhà).'
$newSheet.Range("E7").Value = 'Người độc thân có sức hút mãnh liệt với người khác phái. Dễ thoát ế trong năm nay.'

$newSheet.Range("A8").Value = 'Ngọ'
$newSheet.Range("B8").Value = 'Năm Tuổi (Trực Thái Tuế). "Lửa thử vàng, gian nan thử sức". Áp lực bủa vây nhưng là cơ hội để chứng minh bản lĩnh.'
$newSheet.Range("C8").Value = 'Có nhiều thử thách, cạnh tranh gay gắt tại nơi làm việc. Cần kiên trì, nóng vội là hỏng việc lớn.'
$newSheet.Range("D8").Value = 'Tài chính biến động mạnh. Có cơ hội kiếm tiền nhanh nhưng rủi ro cao. Cẩn thận mất mát đồ đạc.'
$newSheet.Range("E8").Value = 'Tâm tính nóng nảy thất thường dễ làm tổn thương người bên cạnh. Cần học cách kiềm chế cảm xúc.'

$newSheet.Range("A9").Value = 'Mùi'
$newSheet.Range("B9").Value = 'Năm Nhị Hợp (Ngọ - Mùi). Mọi việc hanh thông, suôn sẻ. Được bạn bè, đối tác hỗ trợ nhiệt tình.'
$newSheet.Range("C9").Value = 'Hợp tác làm ăn thuận lợi. Ký kết được nhiều hợp đồng giá trị. Sự nghiệp thăng tiến như diều gặp gió.'
$newSheet.Range("D9").Value = 'Tiền bạc rủng rỉnh. Đầu tư sinh lời tốt. Có thể mua sắm nhà cửa, xe cộ trong năm nay.'
$newSheet.Range("E9").Value = 'Tình cảm thăng hoa. Gia đình hạnh phúc, con cái ngoan ngoãn, thành đạt.'

$newSheet.Range("A10").Value = 'Thân'
$newSheet.Range("B10").Value = 'Một năm Bình Ổn. Sức khỏe và tinh thần đều tốt. Thích hợp để học tập, trau dồi kỹ năng mới.'
$newSheet.Range("C10").Value = 'Công việc không có nhiều biến động. Thích hợp để củng cố vị trí hiện tại hơn là nhảy việc.'
$newSheet.Range("D10").Value = 'Tài lộc ở mức khá. Có lộc ăn uống, quà cáp biếu tặng thường xuyên.'
$newSheet.Range("E10").Value = 'Nhân duyên tốt đẹp. Người độc thân dễ gặp được "nửa kia" qua sự giới thiệu của bạn bè.'

$newSheet.Range("A11").Value = 'Dậu'
$newSheet.Range("B11").Value = 'Năm có sao Đào Hoa chiếu mệnh. Rất lợi cho những người làm nghệ thuật, kinh doanh, ngoại giao.'
$newSheet.Range("C11").Value = 'Được lòng sếp và đồng nghiệp. Công việc trôi chảy nhờ khéo léo trong giao tiếp.'
$newSheet.Range("D11").Value = 'Tài lộc khởi sắc, đặc biệt là các nguồn thu phụ hoặc kinh doanh online.'
$newSheet.Range("E11").Value = 'Đào hoa quá vượng đôi khi lại phiền phức cho người đã kết hôn (dễ vướng thị phi tình ái). Cần giữ mình.'

$newSheet.Range("A12").Value = 'Tuất'
$newSheet.Range("B12").Value = 'Năm Tam Hợp (Dần - Ngọ - Tuất). Vận khí cực thịnh, làm ít hưởng nhiều. Là năm gặt hái thành quả sau bao năm vất vả.'
$newSheet.Range("C12").Value = 'Cơ hội thăng tiến rõ rệt. Có thể được đề bạt lên vị trí lãnh đạo hoặc quản lý cấp cao.'
$newSheet.Range("D12").Value = 'Tiền vào như nước. Có duyên với đất đai, điền sản. Đầu tư đâu thắng đó.'
$newSheet.Range("E12").Value = 'Gia đạo an vui. Có tin vui về chuyện con cái hoặc hỷ sự trong dòng họ.'

$newSheet.Range("A13").Value = 'Hợi'
$newSheet.Range("B13").Value = 'Năm Bình Hòa. Cần chú ý vấn đề sức khỏe và cân bằng cuộc sống - công việc.'
$newSheet.Range("C13").Value = 'Công việc ổn định. Cần tránh va chạm với người có chức quyền. "Dĩ hòa vi quý".'
$newSheet.Range("D13").Value = 'Tài chính đủ chi tiêu, không dư dả nhiều. Tránh cho vay mượn kẻo mất cả tiền lẫn bạn.'
$newSheet.Range("E13").Value = 'Tình cảm bình lặng. Cần quan tâm hơn đến sức khỏe của người lớn tuổi trong nhà.'

# Column widths to roughly match the authored layout
$newSheet.Columns.Item(2).ColumnWidth = 41
$newSheet.Columns.Item(3).ColumnWidth = 30
$newSheet.Range("D1:E1").EntireColumn.ColumnWidth = 8

# Page setup: portrait orientation
$newSheet.PageSetup.Orientation = 1

# Leave the selection on the last written cell (E13), matching the authored file
[void]$newSheet.Range("E13").Select()